# Update "想去人数" (number of people interested) counts for the
# "南宁·烨岑M动漫嘉年华·万圣派对" event on 2024-10-26 (285 -> 287)
# and the "南宁·万圣漫控嘉年华10" event on 2024-11-02 (960 -> 964),
# on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 287
$wsExpo.Range("F4").Value = 964

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 287
$wsAll.Range("F5").Value = 964
